# group smart rule verme tamma
# Update target system addresses on Sheet1 and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 (administrator / bastion) now points at the new sinav host.
$ws.Range("B10").Value = "sinav.quasys.local"

# Row 2 (pam118064) now points at the new internal address.
$ws.Range("B2").Value = "192.168.10.125"

# Move the active selection to C13, matching the saved workbook state.
$ws.Range("C13").Select()
